$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.944.30"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "2.045.97"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.05"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.662"
$ws.Range("E6").Value = "  +0.95%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.44"
$ws.Range("E8").Value = "  +1.57%  "

$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("E10").Value = "  +2.91%  "

$ws.Range("E11").Value = "  +1.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.75"
$ws.Range("E12").Value = "  +4.11%  "

$ws.Range("D13").Value = "2.342.12"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.58"
$ws.Range("E14").Value = "  +6.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.789"
$ws.Range("E15").Value = "  -4.29%  "

$ws.Range("D16").Value = "2.050.24"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").Value = "36.912.52"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.19"
$ws.Range("E18").Value = "  +13.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.02"
$ws.Range("E19").Value = "  +2.18%  "

$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("E21").Value = "  +0.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.89"
$ws.Range("E22").Value = "  -0.81%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("E24").Value = "  -3.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  +9.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.67"
$ws.Range("E26").Value = "  -1.65%  "

$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.78"
$ws.Range("E28").Value = "  -2.52%  "

$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.11"
$ws.Range("E30").Value = "  +5.90%  "

$ws.Range("E31").Value = "  +0.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0612"
$ws.Range("E32").Value = "  -2.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.41"
$ws.Range("E33").Value = "  +1.17%  "

$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0872"
$ws.Range("E35").Value = "  +2.47%  "

$ws.Range("E36").Value = "  -2.99%  "

$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.34"
$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("E39").Value = "  -4.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.17"
$ws.Range("E40").Value = "  +13.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.92"
$ws.Range("E41").Value = "  +24.02%  "

$ws.Range("E42").Value = "  -2.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.22"
$ws.Range("E43").Value = "  -5.76%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.11"
$ws.Range("E44").Value = "  -2.71%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.39"
$ws.Range("E45").Value = "  -1.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("E46").Value = "  +2.29%  "

$ws.Range("D47").Value = "1.278.42"
$ws.Range("E47").Value = "  -1.97%  "

$ws.Range("E48").Value = "  -2.43%  "

$ws.Range("D49").Value = "2.232.72"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  -3.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.51"
$ws.Range("E51").Value = "  -19.83%  "
